$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original genotype labels currently stored in column A (rows 2-15). These get
# moved over to a new "members" column (I), while column A is replaced with
# generic genotype-N identifiers.
$originalLabels = @("1","7","6","14","10","12","5","16","3|2","4|8","21|18","15|11","20|13","17|9|19")

# New header cell I1 = "members", styled like the other header cells in row 1.
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "members"

# Rename the column A identity labels to the generic genotype-N scheme first
# (so the new shared strings "members"/"genotype-N" are interned before the
# original labels get re-inserted under column I, matching insertion order).
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "genotype-$($r - 1)"
}

for ($r = 2; $r -le 15; $r++) {
    $label = $originalLabels[$r - 2]

    # Write the original label into column I as a genuine text value (not a
    # number), without disturbing the cell's style: build it via a formula
    # that evaluates to text, then convert the formula to a static value via
    # copy / paste-special.
    $cellI = $ws.Cells.Item($r, 9)
    $cellI.Formula = '="' + $label + '"'
    $cellI.Copy()
    $cellI.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
